$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 393-407 ---
$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A393:A407").NumberFormat = "@"
$ws.Cells.Item(393, 1).Value = "2026-02-06"
$ws.Cells.Item(393, 2).Value = "10:12:37"
$ws.Cells.Item(393, 3).Value = "10:00"
$ws.Cells.Item(393, 4).Value = "Bathroom"
$ws.Cells.Item(393, 5).Value = "No Motion"
$ws.Cells.Item(393, 6).Value = "Inactive"
$ws.Cells.Item(394, 1).Value = "2026-02-06"
$ws.Cells.Item(394, 2).Value = "10:12:39"
$ws.Cells.Item(394, 3).Value = "10:00"
$ws.Cells.Item(394, 4).Value = "Bathroom"
$ws.Cells.Item(394, 5).Value = "No Motion"
$ws.Cells.Item(394, 6).Value = "Inactive"
$ws.Cells.Item(395, 1).Value = "2026-02-06"
$ws.Cells.Item(395, 2).Value = "10:12:40"
$ws.Cells.Item(395, 3).Value = "10:00"
$ws.Cells.Item(395, 4).Value = "Bathroom"
$ws.Cells.Item(395, 5).Value = "No Motion"
$ws.Cells.Item(395, 6).Value = "Inactive"
$ws.Cells.Item(396, 1).Value = "2026-02-06"
$ws.Cells.Item(396, 2).Value = "10:12:45"
$ws.Cells.Item(396, 3).Value = "10:00"
$ws.Cells.Item(396, 4).Value = "Bathroom"
$ws.Cells.Item(396, 5).Value = "No Motion"
$ws.Cells.Item(396, 6).Value = "Inactive"
$ws.Cells.Item(397, 1).Value = "2026-02-06"
$ws.Cells.Item(397, 2).Value = "10:12:50"
$ws.Cells.Item(397, 3).Value = "10:00"
$ws.Cells.Item(397, 4).Value = "Bathroom"
$ws.Cells.Item(397, 5).Value = "No Motion"
$ws.Cells.Item(397, 6).Value = "Inactive"
$ws.Cells.Item(398, 1).Value = "2026-02-06"
$ws.Cells.Item(398, 2).Value = "10:12:56"
$ws.Cells.Item(398, 3).Value = "10:00"
$ws.Cells.Item(398, 4).Value = "Bathroom"
$ws.Cells.Item(398, 5).Value = "No Motion"
$ws.Cells.Item(398, 6).Value = "Inactive"
$ws.Cells.Item(399, 1).Value = "2026-02-06"
$ws.Cells.Item(399, 2).Value = "10:13:00"
$ws.Cells.Item(399, 3).Value = "10:00"
$ws.Cells.Item(399, 4).Value = "Bathroom"
$ws.Cells.Item(399, 5).Value = "No Motion"
$ws.Cells.Item(399, 6).Value = "Inactive"
$ws.Cells.Item(400, 1).Value = "2026-02-06"
$ws.Cells.Item(400, 2).Value = "10:13:01"
$ws.Cells.Item(400, 3).Value = "10:00"
$ws.Cells.Item(400, 4).Value = "Bathroom"
$ws.Cells.Item(400, 5).Value = "Motion Detected"
$ws.Cells.Item(400, 6).Value = "Active"
$ws.Cells.Item(401, 1).Value = "2026-02-06"
$ws.Cells.Item(401, 2).Value = "10:13:08"
$ws.Cells.Item(401, 3).Value = "10:00"
$ws.Cells.Item(401, 4).Value = "Bathroom"
$ws.Cells.Item(401, 5).Value = "No Motion"
$ws.Cells.Item(401, 6).Value = "Inactive"
$ws.Cells.Item(402, 1).Value = "2026-02-06"
$ws.Cells.Item(402, 2).Value = "10:13:13"
$ws.Cells.Item(402, 3).Value = "10:00"
$ws.Cells.Item(402, 4).Value = "Bathroom"
$ws.Cells.Item(402, 5).Value = "No Motion"
$ws.Cells.Item(402, 6).Value = "Inactive"
$ws.Cells.Item(403, 1).Value = "2026-02-06"
$ws.Cells.Item(403, 2).Value = "10:13:20"
$ws.Cells.Item(403, 3).Value = "10:00"
$ws.Cells.Item(403, 4).Value = "Bathroom"
$ws.Cells.Item(403, 5).Value = "No Motion"
$ws.Cells.Item(403, 6).Value = "Inactive"
$ws.Cells.Item(404, 1).Value = "2026-02-06"
$ws.Cells.Item(404, 2).Value = "10:13:21"
$ws.Cells.Item(404, 3).Value = "10:00"
$ws.Cells.Item(404, 4).Value = "Bathroom"
$ws.Cells.Item(404, 5).Value = "Motion Detected"
$ws.Cells.Item(404, 6).Value = "Active"
$ws.Cells.Item(405, 1).Value = "2026-02-06"
$ws.Cells.Item(405, 2).Value = "10:13:28"
$ws.Cells.Item(405, 3).Value = "10:00"
$ws.Cells.Item(405, 4).Value = "Bathroom"
$ws.Cells.Item(405, 5).Value = "No Motion"
$ws.Cells.Item(405, 6).Value = "Inactive"
$ws.Cells.Item(406, 1).Value = "2026-02-06"
$ws.Cells.Item(406, 2).Value = "10:13:33"
$ws.Cells.Item(406, 3).Value = "10:00"
$ws.Cells.Item(406, 4).Value = "Bathroom"
$ws.Cells.Item(406, 5).Value = "No Motion"
$ws.Cells.Item(406, 6).Value = "Inactive"
$ws.Cells.Item(407, 1).Value = "2026-02-06"
$ws.Cells.Item(407, 2).Value = "10:13:35"
$ws.Cells.Item(407, 3).Value = "10:00"
$ws.Cells.Item(407, 4).Value = "Bathroom"
$ws.Cells.Item(407, 5).Value = "Motion Detected"
$ws.Cells.Item(407, 6).Value = "Active"

# --- Humidity sheet: append rows 271-277 ---
$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A271:A277").NumberFormat = "@"
$ws.Range("E271").NumberFormat = "@"
$ws.Range("E272").NumberFormat = "@"
$ws.Range("E273").NumberFormat = "@"
$ws.Range("E274").NumberFormat = "@"
$ws.Range("E275").NumberFormat = "@"
$ws.Range("E276").NumberFormat = "@"
$ws.Range("E277").NumberFormat = "@"
$ws.Cells.Item(271, 1).Value = "2026-02-06"
$ws.Cells.Item(271, 2).Value = "10:12:38"
$ws.Cells.Item(271, 3).Value = "10:00"
$ws.Cells.Item(271, 4).Value = "Bathroom"
$ws.Cells.Item(271, 5).Value = "68.0%"
$ws.Cells.Item(271, 6).Value = "Active"
$ws.Cells.Item(272, 1).Value = "2026-02-06"
$ws.Cells.Item(272, 2).Value = "10:12:43"
$ws.Cells.Item(272, 3).Value = "10:00"
$ws.Cells.Item(272, 4).Value = "Bathroom"
$ws.Cells.Item(272, 5).Value = "68.0%"
$ws.Cells.Item(272, 6).Value = "Active"
$ws.Cells.Item(273, 1).Value = "2026-02-06"
$ws.Cells.Item(273, 2).Value = "10:12:48"
$ws.Cells.Item(273, 3).Value = "10:00"
$ws.Cells.Item(273, 4).Value = "Bathroom"
$ws.Cells.Item(273, 5).Value = "67.5%"
$ws.Cells.Item(273, 6).Value = "Active"
$ws.Cells.Item(274, 1).Value = "2026-02-06"
$ws.Cells.Item(274, 2).Value = "10:12:58"
$ws.Cells.Item(274, 3).Value = "10:00"
$ws.Cells.Item(274, 4).Value = "Bathroom"
$ws.Cells.Item(274, 5).Value = "67.7%"
$ws.Cells.Item(274, 6).Value = "Active"
$ws.Cells.Item(275, 1).Value = "2026-02-06"
$ws.Cells.Item(275, 2).Value = "10:13:18"
$ws.Cells.Item(275, 3).Value = "10:00"
$ws.Cells.Item(275, 4).Value = "Bathroom"
$ws.Cells.Item(275, 5).Value = "69.5%"
$ws.Cells.Item(275, 6).Value = "Active"
$ws.Cells.Item(276, 1).Value = "2026-02-06"
$ws.Cells.Item(276, 2).Value = "10:13:29"
$ws.Cells.Item(276, 3).Value = "10:00"
$ws.Cells.Item(276, 4).Value = "Bathroom"
$ws.Cells.Item(276, 5).Value = "69.7%"
$ws.Cells.Item(276, 6).Value = "Active"
$ws.Cells.Item(277, 1).Value = "2026-02-06"
$ws.Cells.Item(277, 2).Value = "10:13:34"
$ws.Cells.Item(277, 3).Value = "10:00"
$ws.Cells.Item(277, 4).Value = "Bathroom"
$ws.Cells.Item(277, 5).Value = "69.6%"
$ws.Cells.Item(277, 6).Value = "Active"

# --- Temperature sheet: append rows 271-277 ---
$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A271:A277").NumberFormat = "@"
$ws.Cells.Item(271, 1).Value = "2026-02-06"
$ws.Cells.Item(271, 2).Value = "10:12:39"
$ws.Cells.Item(271, 3).Value = "10:00"
$ws.Cells.Item(271, 4).Value = "Bathroom"
$ws.Cells.Item(271, 5).Value = "28.0C"
$ws.Cells.Item(271, 6).Value = "Active"
$ws.Cells.Item(272, 1).Value = "2026-02-06"
$ws.Cells.Item(272, 2).Value = "10:12:44"
$ws.Cells.Item(272, 3).Value = "10:00"
$ws.Cells.Item(272, 4).Value = "Bathroom"
$ws.Cells.Item(272, 5).Value = "28.0C"
$ws.Cells.Item(272, 6).Value = "Active"
$ws.Cells.Item(273, 1).Value = "2026-02-06"
$ws.Cells.Item(273, 2).Value = "10:12:49"
$ws.Cells.Item(273, 3).Value = "10:00"
$ws.Cells.Item(273, 4).Value = "Bathroom"
$ws.Cells.Item(273, 5).Value = "28.0C"
$ws.Cells.Item(273, 6).Value = "Active"
$ws.Cells.Item(274, 1).Value = "2026-02-06"
$ws.Cells.Item(274, 2).Value = "10:12:59"
$ws.Cells.Item(274, 3).Value = "10:00"
$ws.Cells.Item(274, 4).Value = "Bathroom"
$ws.Cells.Item(274, 5).Value = "28.0C"
$ws.Cells.Item(274, 6).Value = "Active"
$ws.Cells.Item(275, 1).Value = "2026-02-06"
$ws.Cells.Item(275, 2).Value = "10:13:19"
$ws.Cells.Item(275, 3).Value = "10:00"
$ws.Cells.Item(275, 4).Value = "Bathroom"
$ws.Cells.Item(275, 5).Value = "28.1C"
$ws.Cells.Item(275, 6).Value = "Active"
$ws.Cells.Item(276, 1).Value = "2026-02-06"
$ws.Cells.Item(276, 2).Value = "10:13:30"
$ws.Cells.Item(276, 3).Value = "10:00"
$ws.Cells.Item(276, 4).Value = "Bathroom"
$ws.Cells.Item(276, 5).Value = "28.1C"
$ws.Cells.Item(276, 6).Value = "Active"
$ws.Cells.Item(277, 1).Value = "2026-02-06"
$ws.Cells.Item(277, 2).Value = "10:13:35"
$ws.Cells.Item(277, 3).Value = "10:00"
$ws.Cells.Item(277, 4).Value = "Bathroom"
$ws.Cells.Item(277, 5).Value = "28.0C"
$ws.Cells.Item(277, 6).Value = "Active"
